# Glossary.xlsx: move the "Use (X) - AllStar prediction" marks in column C
# to the correct rows, and update the active selection to C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# The "X" mark was on the wrong row for a few glossary entries;
# clear it from the old (incorrect) rows...
$ws.Range("C18").Clear()
$ws.Range("C38").Clear()

# ...and set it on the correct rows.
$ws.Range("C28").Value = "X"
$ws.Range("C47").Value = "X"
$ws.Range("C70").Value = "X"

# Update the current selection to match the saved workbook state (C1).
$ws.Range("C1").Select()
